# Sync automatico del tracker (cada 3h)
# Appends the latest scraped match rows (event_id, fecha, jugador_A,
# jugador_B, pronostico, cuota) to the bottom of the tracker sheet.
# resultado/profit (G/H) are left blank until the match finishes.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("14552518", "2025-09-04", "Tsung-Hao Huang",   "Bernard Tomic",  "Gana Tsung-Hao Huang",   3.75),
    @("14552529", "2025-09-04", "Daniel Evans",      "Filip Peliwo",   "Gana Filip Peliwo",      5.5),
    @("14552909", "2025-09-04", "Tom Paris",         "Mark Lajal",     "Gana Tom Paris",         3.5),
    @("14552660", "2025-09-04", "Mert Naci Turker",  "Calvin Hemery",  "Gana Mert Naci Turker",  5.5)
)

$startRow = 121

for ($i = 0; $i -lt $rows.Count; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]

    # event_id / fecha / jugador_A / jugador_B / pronostico must stay as
    # plain text (event_id and fecha look numeric/date-like, so force the
    # "Text" number format before writing them, then drop back to the
    # Normal style so no extra formatting is left behind on the cell).
    $textRange = $ws.Range("A" + $r + ":E" + $r)
    $textRange.NumberFormat = "@"

    $ws.Range("A" + $r).Value = $data[0]
    $ws.Range("B" + $r).Value = $data[1]
    $ws.Range("C" + $r).Value = $data[2]
    $ws.Range("D" + $r).Value = $data[3]
    $ws.Range("E" + $r).Value = $data[4]

    $textRange.Style = "Normal"

    # cuota (F) is numeric.
    $ws.Range("F" + $r).Value = $data[5]

    # resultado / profit (G/H): match still pending, no value yet.
}
